$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally carried 5 pairs of "외주" (outsourced) price/comparison
# columns (단가/대비%) in L:U. This edit drops the "대비%" (comparison) column of
# each pair, keeping only the "단가" (price) column, and renames the remaining
# price-column headers from "[단가(외주N)" / "대비(외주N)]" bracket pairs to the
# plain "단가(외주N)" label. Net effect: columns M, O, Q, S, U are removed and
# L:P become "단가(외주1)".."단가(외주5)".
#
# Delete right-to-left so earlier deletes don't invalidate later column letters.
$ws.Columns("U").Delete()
$ws.Columns("S").Delete()
$ws.Columns("Q").Delete()
$ws.Columns("O").Delete()
$ws.Columns("M").Delete()

# Re-label the surviving price-column headers (row 1, columns L:P).
$ws.Range("L1").Value2 = "단가(외주1)"
$ws.Range("M1").Value2 = "단가(외주2)"
$ws.Range("N1").Value2 = "단가(외주3)"
$ws.Range("O1").Value2 = "단가(외주4)"
$ws.Range("P1").Value2 = "단가(외주5)"

# Column J ("대비(표준,%)") was widened to fit its content.
$ws.Columns("J").ColumnWidth = 14.318181818181818

# Selection cursor moved to K12.
$ws.Range("K12").Select()
